$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 13.458797
$ws.Cells.Item(2, 8).Value = 40.376391
$ws.Cells.Item(2, 9).Value = 0.06830096976102129
$ws.Cells.Item(2, 10).Value = 0.06973720484213804
$ws.Cells.Item(2, 13).Value = 28.85518433333334
$ws.Cells.Item(2, 14).Value = 86.56555300000001
$ws.Cells.Item(2, 15).Value = 0.1999651185353207
$ws.Cells.Item(2, 16).Value = 0.2044513327926365
$ws.Cells.Item(2, 17).Value = 388.3560683399137
$ws.Cells.Item(2, 18).Value = 3495.204615059223
$ws.Cells.Item(2, 19).Value = 0.01365781151433998
$ws.Cells.Item(2, 20).Value = 0.01425786447520823

$ws.Cells.Item(3, 7).Value = 13.458797
$ws.Cells.Item(3, 8).Value = 40.376391
$ws.Cells.Item(3, 9).Value = 0.06830096976102129
$ws.Cells.Item(3, 10).Value = 0.06973720484213804
$ws.Cells.Item(3, 15).Value = 0.3546352265743414
$ws.Cells.Item(3, 16).Value = 0.3625914622481308
$ws.Cells.Item(3, 17).Value = 688.7438334047187
$ws.Cells.Item(3, 18).Value = 6198.694500642468
$ws.Cells.Item(3, 19).Value = 0.02422192988644703
$ws.Cells.Item(3, 20).Value = 0.02528611507680826

$ws.Cells.Item(4, 7).Value = 13.458797
$ws.Cells.Item(4, 8).Value = 40.376391
$ws.Cells.Item(4, 9).Value = 0.06830096976102129
$ws.Cells.Item(4, 10).Value = 0.06973720484213804
$ws.Cells.Item(4, 13).Value = 29.393479
$ws.Cells.Item(4, 14).Value = 88.180437
$ws.Cells.Item(4, 15).Value = 0.2036954761578358
$ws.Cells.Item(4, 16).Value = 0.2082653809291453
$ws.Cells.Item(4, 17).Value = 395.6008669847629
$ws.Cells.Item(4, 18).Value = 3560.407802862867
$ws.Cells.Item(4, 19).Value = 0.01391259855751318
$ws.Cells.Item(4, 20).Value = 0.01452384553138171

$ws.Cells.Item(5, 7).Value = 13.458797
$ws.Cells.Item(5, 8).Value = 40.376391
$ws.Cells.Item(5, 9).Value = 0.06830096976102129
$ws.Cells.Item(5, 10).Value = 0.06973720484213804
$ws.Cells.Item(5, 13).Value = 9.499066500000001
$ws.Cells.Item(5, 14).Value = 18.998133
$ws.Cells.Item(5, 15).Value = 0.0658280999596015
$ws.Cells.Item(5, 16).Value = 0.04486996822421697
$ws.Cells.Item(5, 17).Value = 127.8460077130005
$ws.Cells.Item(5, 18).Value = 767.0760462780031
$ws.Cells.Item(5, 19).Value = 0.004496123064766228
$ws.Cells.Item(5, 20).Value = 0.003129106165312443

$ws.Cells.Item(6, 7).Value = 13.458797
$ws.Cells.Item(6, 8).Value = 40.376391
$ws.Cells.Item(6, 9).Value = 0.06830096976102129
$ws.Cells.Item(6, 10).Value = 0.06973720484213804
$ws.Cells.Item(6, 13).Value = 25.37910966666666
$ws.Cells.Item(6, 14).Value = 76.13732899999999
$ws.Cells.Item(6, 15).Value = 0.1758760787729007
$ws.Cells.Item(6, 16).Value = 0.1798218558058706
$ws.Cells.Item(6, 17).Value = 341.5722850444043
$ws.Cells.Item(6, 18).Value = 3074.150565399639
$ws.Cells.Item(6, 19).Value = 0.01201250673795489
$ws.Cells.Item(6, 20).Value = 0.01254027359342741

$ws.Cells.Item(7, 9).Value = 0.1240039124627887
$ws.Cells.Item(7, 10).Value = 0.1266114708898203
$ws.Cells.Item(7, 13).Value = 28.85518433333334
$ws.Cells.Item(7, 14).Value = 86.56555300000001
$ws.Cells.Item(7, 15).Value = 0.1999651185353207
$ws.Cells.Item(7, 16).Value = 0.2044513327926365
$ws.Cells.Item(7, 17).Value = 705.0803534900695
$ws.Cells.Item(7, 18).Value = 6345.723181410624
$ws.Cells.Item(7, 19).Value = 0.02479645705446507
$ws.Cells.Item(7, 20).Value = 0.02588588397025986

$ws.Cells.Item(8, 9).Value = 0.1240039124627887
$ws.Cells.Item(8, 10).Value = 0.1266114708898203
$ws.Cells.Item(8, 15).Value = 0.3546352265743414
$ws.Cells.Item(8, 16).Value = 0.3625914622481308
$ws.Cells.Item(8, 17).Value = 1250.449742157909
$ws.Cells.Item(8, 18).Value = 11254.04767942118
$ws.Cells.Item(8, 19).Value = 0.04397615559234587
$ws.Cells.Item(8, 20).Value = 0.0459082383673266

$ws.Cells.Item(9, 9).Value = 0.1240039124627887
$ws.Cells.Item(9, 10).Value = 0.1266114708898203
$ws.Cells.Item(9, 13).Value = 29.393479
$ws.Cells.Item(9, 14).Value = 88.180437
$ws.Cells.Item(9, 15).Value = 0.2036954761578358
$ws.Cells.Item(9, 16).Value = 0.2082653809291453
$ws.Cells.Item(9, 17).Value = 718.233656878144
$ws.Cells.Item(9, 18).Value = 6464.102911903296
$ws.Cells.Item(9, 19).Value = 0.02525903599454234
$ws.Cells.Item(9, 20).Value = 0.02636878621486782

$ws.Cells.Item(10, 9).Value = 0.1240039124627887
$ws.Cells.Item(10, 10).Value = 0.1266114708898203
$ws.Cells.Item(10, 13).Value = 9.499066500000001
$ws.Cells.Item(10, 14).Value = 18.998133
$ws.Cells.Item(10, 15).Value = 0.0658280999596015
$ws.Cells.Item(10, 16).Value = 0.04486996822421697
$ws.Cells.Item(10, 17).Value = 232.110981800544
$ws.Cells.Item(10, 18).Value = 1392.665890803264
$ws.Cells.Item(10, 19).Value = 0.008162941944982127
$ws.Cells.Item(10, 20).Value = 0.00568105267564761

$ws.Cells.Item(11, 9).Value = 0.1240039124627887
$ws.Cells.Item(11, 10).Value = 0.1266114708898203
$ws.Cells.Item(11, 13).Value = 25.37910966666666
$ws.Cells.Item(11, 14).Value = 76.13732899999999
$ws.Cells.Item(11, 15).Value = 0.1758760787729007
$ws.Cells.Item(11, 16).Value = 0.1798218558058706
$ws.Cells.Item(11, 17).Value = 620.1419962639146
$ws.Cells.Item(11, 18).Value = 5581.277966375232
$ws.Cells.Item(11, 19).Value = 0.0218093218764533
$ws.Cells.Item(11, 20).Value = 0.02276750966171846

$ws.Cells.Item(12, 7).Value = 76.51423666666666
$ws.Cells.Item(12, 8).Value = 229.54271
$ws.Cells.Item(12, 9).Value = 0.3882959647030583
$ws.Cells.Item(12, 10).Value = 0.3964610652618627
$ws.Cells.Item(12, 13).Value = 28.85518433333334
$ws.Cells.Item(12, 14).Value = 86.56555300000001
$ws.Cells.Item(12, 15).Value = 0.1999651185353207
$ws.Cells.Item(12, 16).Value = 0.2044513327926365
$ws.Cells.Item(12, 17).Value = 2207.832403140959
$ws.Cells.Item(12, 18).Value = 19870.49162826863
$ws.Cells.Item(12, 19).Value = 0.07764564860863375
$ws.Cells.Item(12, 20).Value = 0.08105699319317627

$ws.Cells.Item(13, 7).Value = 76.51423666666666
$ws.Cells.Item(13, 8).Value = 229.54271
$ws.Cells.Item(13, 9).Value = 0.3882959647030583
$ws.Cells.Item(13, 10).Value = 0.3964610652618627
$ws.Cells.Item(13, 15).Value = 0.3546352265743414
$ws.Cells.Item(13, 16).Value = 0.3625914622481308
$ws.Cells.Item(13, 17).Value = 3915.558624729676
$ws.Cells.Item(13, 18).Value = 35240.02762256708
$ws.Cells.Item(13, 19).Value = 0.1377034274203716
$ws.Cells.Item(13, 20).Value = 0.1437533973777504

$ws.Cells.Item(14, 7).Value = 76.51423666666666
$ws.Cells.Item(14, 8).Value = 229.54271
$ws.Cells.Item(14, 9).Value = 0.3882959647030583
$ws.Cells.Item(14, 10).Value = 0.3964610652618627
$ws.Cells.Item(14, 13).Value = 29.393479
$ws.Cells.Item(14, 14).Value = 88.180437
$ws.Cells.Item(14, 15).Value = 0.2036954761578358
$ws.Cells.Item(14, 16).Value = 0.2082653809291453
$ws.Cells.Item(14, 17).Value = 2249.019608662697
$ws.Cells.Item(14, 18).Value = 20241.17647796427
$ws.Cells.Item(14, 19).Value = 0.07909413142035568
$ws.Cells.Item(14, 20).Value = 0.08256911478033656

$ws.Cells.Item(15, 7).Value = 76.51423666666666
$ws.Cells.Item(15, 8).Value = 229.54271
$ws.Cells.Item(15, 9).Value = 0.3882959647030583
$ws.Cells.Item(15, 10).Value = 0.3964610652618627
$ws.Cells.Item(15, 13).Value = 9.499066500000001
$ws.Cells.Item(15, 14).Value = 18.998133
$ws.Cells.Item(15, 15).Value = 0.0658280999596015
$ws.Cells.Item(15, 16).Value = 0.04486996822421697
$ws.Cells.Item(15, 17).Value = 726.8138222934051
$ws.Cells.Item(15, 18).Value = 4360.88293376043
$ws.Cells.Item(15, 19).Value = 0.02556078557838282
$ws.Cells.Item(15, 20).Value = 0.01778919540043899

$ws.Cells.Item(16, 7).Value = 76.51423666666666
$ws.Cells.Item(16, 8).Value = 229.54271
$ws.Cells.Item(16, 9).Value = 0.3882959647030583
$ws.Cells.Item(16, 10).Value = 0.3964610652618627
$ws.Cells.Item(16, 13).Value = 25.37910966666666
$ws.Cells.Item(16, 14).Value = 76.13732899999999
$ws.Cells.Item(16, 15).Value = 0.1758760787729007
$ws.Cells.Item(16, 16).Value = 0.1798218558058706
$ws.Cells.Item(16, 17).Value = 1941.863203424621
$ws.Cells.Item(16, 18).Value = 17476.76883082159
$ws.Cells.Item(16, 19).Value = 0.06829197167531453
$ws.Cells.Item(16, 20).Value = 0.07129236451016054

$ws.Cells.Item(17, 7).Value = 12.174794
$ws.Cells.Item(17, 8).Value = 24.349588
$ws.Cells.Item(17, 9).Value = 0.06178488588843889
$ws.Cells.Item(17, 10).Value = 0.04205606702633888
$ws.Cells.Item(17, 13).Value = 28.85518433333334
$ws.Cells.Item(17, 14).Value = 86.56555300000001
$ws.Cells.Item(17, 15).Value = 0.1999651185353207
$ws.Cells.Item(17, 16).Value = 0.2044513327926365
$ws.Cells.Item(17, 17).Value = 351.3059250903607
$ws.Cells.Item(17, 18).Value = 2107.835550542164
$ws.Cells.Item(17, 19).Value = 0.01235482203037294
$ws.Cells.Item(17, 20).Value = 0.008598418955551438

$ws.Cells.Item(18, 7).Value = 12.174794
$ws.Cells.Item(18, 8).Value = 24.349588
$ws.Cells.Item(18, 9).Value = 0.06178488588843889
$ws.Cells.Item(18, 10).Value = 0.04205606702633888
$ws.Cells.Item(18, 15).Value = 0.3546352265743414
$ws.Cells.Item(18, 16).Value = 0.3625914622481308
$ws.Cells.Item(18, 17).Value = 623.0359437379707
$ws.Cells.Item(18, 18).Value = 3738.215662427824
$ws.Cells.Item(18, 19).Value = 0.02191109700591636
$ws.Cells.Item(18, 20).Value = 0.01524917083948562

$ws.Cells.Item(19, 7).Value = 12.174794
$ws.Cells.Item(19, 8).Value = 24.349588
$ws.Cells.Item(19, 9).Value = 0.06178488588843889
$ws.Cells.Item(19, 10).Value = 0.04205606702633888
$ws.Cells.Item(19, 13).Value = 29.393479
$ws.Cells.Item(19, 14).Value = 88.180437
$ws.Cells.Item(19, 15).Value = 0.2036954761578358
$ws.Cells.Item(19, 16).Value = 0.2082653809291453
$ws.Cells.Item(19, 17).Value = 357.859551768326
$ws.Cells.Item(19, 18).Value = 2147.157310609956
$ws.Cells.Item(19, 19).Value = 0.01258530175040311
$ws.Cells.Item(19, 20).Value = 0.008758822819622133

$ws.Cells.Item(20, 7).Value = 12.174794
$ws.Cells.Item(20, 8).Value = 24.349588
$ws.Cells.Item(20, 9).Value = 0.06178488588843889
$ws.Cells.Item(20, 10).Value = 0.04205606702633888
$ws.Cells.Item(20, 13).Value = 9.499066500000001
$ws.Cells.Item(20, 14).Value = 18.998133
$ws.Cells.Item(20, 15).Value = 0.0658280999596015
$ws.Cells.Item(20, 16).Value = 0.04486996822421697
$ws.Cells.Item(20, 17).Value = 115.649177829801
$ws.Cells.Item(20, 18).Value = 462.5967113192041
$ws.Cells.Item(20, 19).Value = 0.004067181644256727
$ws.Cells.Item(20, 20).Value = 0.001887054391107365

$ws.Cells.Item(21, 7).Value = 12.174794
$ws.Cells.Item(21, 8).Value = 24.349588
$ws.Cells.Item(21, 9).Value = 0.06178488588843889
$ws.Cells.Item(21, 10).Value = 0.04205606702633888
$ws.Cells.Item(21, 13).Value = 25.37910966666666
$ws.Cells.Item(21, 14).Value = 76.13732899999999
$ws.Cells.Item(21, 15).Value = 0.1758760787729007
$ws.Cells.Item(21, 16).Value = 0.1798218558058706
$ws.Cells.Item(21, 17).Value = 308.9854320950753
$ws.Cells.Item(21, 18).Value = 1853.912592570452
$ws.Cells.Item(21, 19).Value = 0.01086648345748976
$ws.Cells.Item(21, 20).Value = 0.007562600020572341

$ws.Cells.Item(22, 7).Value = 70.46836733333333
$ws.Cells.Item(22, 8).Value = 211.405102
$ws.Cells.Item(22, 9).Value = 0.3576142671846927
$ws.Cells.Item(22, 10).Value = 0.36513419197984
$ws.Cells.Item(22, 13).Value = 28.85518433333334
$ws.Cells.Item(22, 14).Value = 86.56555300000001
$ws.Cells.Item(22, 15).Value = 0.1999651185353207
$ws.Cells.Item(22, 16).Value = 0.2044513327926365
$ws.Cells.Item(22, 17).Value = 2033.377729072379
$ws.Cells.Item(22, 18).Value = 18300.39956165141
$ws.Cells.Item(22, 19).Value = 0.07151037932750892
$ws.Cells.Item(22, 20).Value = 0.0746521721984407

$ws.Cells.Item(23, 7).Value = 70.46836733333333
$ws.Cells.Item(23, 8).Value = 211.405102
$ws.Cells.Item(23, 9).Value = 0.3576142671846927
$ws.Cells.Item(23, 10).Value = 0.36513419197984
$ws.Cells.Item(23, 15).Value = 0.3546352265743414
$ws.Cells.Item(23, 16).Value = 0.3625914622481308
$ws.Cells.Item(23, 17).Value = 3606.165800028922
$ws.Cells.Item(23, 18).Value = 32455.4922002603
$ws.Cells.Item(23, 19).Value = 0.1268226166692606
$ws.Cells.Item(23, 20).Value = 0.1323945405867599

$ws.Cells.Item(24, 7).Value = 70.46836733333333
$ws.Cells.Item(24, 8).Value = 211.405102
$ws.Cells.Item(24, 9).Value = 0.3576142671846927
$ws.Cells.Item(24, 10).Value = 0.36513419197984
$ws.Cells.Item(24, 13).Value = 29.393479
$ws.Cells.Item(24, 14).Value = 88.180437
$ws.Cells.Item(24, 15).Value = 0.2036954761578358
$ws.Cells.Item(24, 16).Value = 0.2082653809291453
$ws.Cells.Item(24, 17).Value = 2071.310475376619
$ws.Cells.Item(24, 18).Value = 18641.79427838957
$ws.Cells.Item(24, 19).Value = 0.07284440843502152
$ws.Cells.Item(24, 20).Value = 0.07604481158293704

$ws.Cells.Item(25, 7).Value = 70.46836733333333
$ws.Cells.Item(25, 8).Value = 211.405102
$ws.Cells.Item(25, 9).Value = 0.3576142671846927
$ws.Cells.Item(25, 10).Value = 0.36513419197984
$ws.Cells.Item(25, 13).Value = 9.499066500000001
$ws.Cells.Item(25, 14).Value = 18.998133
$ws.Cells.Item(25, 15).Value = 0.0658280999596015
$ws.Cells.Item(25, 16).Value = 0.04486996822421697
$ws.Cells.Item(25, 17).Value = 669.383707445761
$ws.Cells.Item(25, 18).Value = 4016.302244674567
$ws.Cells.Item(25, 19).Value = 0.02354106772721359
$ws.Cells.Item(25, 20).Value = 0.01638355959171056

$ws.Cells.Item(26, 7).Value = 70.46836733333333
$ws.Cells.Item(26, 8).Value = 211.405102
$ws.Cells.Item(26, 9).Value = 0.3576142671846927
$ws.Cells.Item(26, 10).Value = 0.36513419197984
$ws.Cells.Item(26, 13).Value = 25.37910966666666
$ws.Cells.Item(26, 14).Value = 76.13732899999999
$ws.Cells.Item(26, 15).Value = 0.1758760787729007
$ws.Cells.Item(26, 16).Value = 0.1798218558058706
$ws.Cells.Item(26, 17).Value = 1788.424422583617
$ws.Cells.Item(26, 18).Value = 16095.81980325256
$ws.Cells.Item(26, 19).Value = 0.06289579502568816
$ws.Cells.Item(26, 20).Value = 0.06565910801999188
